# Fruta / hortaliza, semanal
# Inserts three new weekly price observations for "Perejil" (Feria Lagunitas de
# Puerto Montt) into the dataset: one near the top of the date range, one in
# the middle, and one appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param($Row, $D, $J, $K, $L, $M, $N, $O, $P, $Q)

    $ws.Cells.Item($Row, 1).Value = 4
    $ws.Cells.Item($Row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($Row, 3).Value = "Los Lagos"
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 10
    $ws.Cells.Item($Row, 6).Value = 100112044
    $ws.Cells.Item($Row, 7).Value = "Perejil"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# 1) Insert a new row at row 48 (pushes the existing rows 48-141 down to 49-142)
$ws.Rows.Item(48).Insert()
Set-DataRow 48 44428 180 4500 4500 4500 "`$/docena de atados (3 kilos)" "Región Metropolitana" 1500 3

# 2) Insert a new row at row 87 (pushes rows 87-142 down to 88-143)
$ws.Rows.Item(87).Insert()
Set-DataRow 87 44435 340 4500 5000 4765 "`$/docena de atados (3 kilos)" "Región Metropolitana" 1588 3

# 3) Append a new row at the end (row 144)
Set-DataRow 144 44432 180 5000 5000 5000 "`$/docena de atados (3 kilos)" "Región Metropolitana" 1667 3
